$d = $word.ActiveDocument
$tbl = $d.Tables.Item(1)

# Word table rows/columns are 1-based. Only the 5 data rows (1, 5, 10, 15, 20)
# contain text; the rows between them are blank spacer rows and are untouched.

# Row 1: 29x31=899, 81x99=8019, 84x95=7980, 43x33=1419, 42x98=4116
$tbl.Cell(1,1).Range.Text = "49×71=3479"
$tbl.Cell(1,2).Range.Text = "71×70=4970"
$tbl.Cell(1,3).Range.Text = "36×55=1980"
$tbl.Cell(1,4).Range.Text = "19×63=1197"
$tbl.Cell(1,5).Range.Text = "23×79=1817"

# Row 5: 90x90=8100, 98x74=7252, 96x49=4704, 41x98=4018, 94x83=7802
$tbl.Cell(5,1).Range.Text = "65×17=1105"
$tbl.Cell(5,2).Range.Text = "36×86=3096"
$tbl.Cell(5,3).Range.Text = "68×85=5780"
$tbl.Cell(5,4).Range.Text = "83×89=7387"
$tbl.Cell(5,5).Range.Text = "21×70=1470"

# Row 10: 34x35=1190, 51x43=2193, 76x80=6080, 66x79=5214, 17x57=969
$tbl.Cell(10,1).Range.Text = "65×13=845"
$tbl.Cell(10,2).Range.Text = "46×29=1334"
$tbl.Cell(10,3).Range.Text = "80×30=2400"
$tbl.Cell(10,4).Range.Text = "17×57=969"
$tbl.Cell(10,5).Range.Text = "21×61=1281"

# Row 15: 71x78=5538, 84x52=4368, 19x94=1786, 34x43=1462, 24x19=456
$tbl.Cell(15,1).Range.Text = "19×20=380"
$tbl.Cell(15,2).Range.Text = "23×32=736"
$tbl.Cell(15,3).Range.Text = "34×67=2278"
$tbl.Cell(15,4).Range.Text = "33×30=990"
$tbl.Cell(15,5).Range.Text = "92×57=5244"

# Row 20: 49x76=3724, 31x91=2821, 25x39=975, 80x19=1520, 95x25=2375
$tbl.Cell(20,1).Range.Text = "97×32=3104"
$tbl.Cell(20,2).Range.Text = "95×55=5225"
$tbl.Cell(20,3).Range.Text = "41×58=2378"
$tbl.Cell(20,4).Range.Text = "83×24=1992"
$tbl.Cell(20,5).Range.Text = "15×53=795"
